# CharacterAbilityDataSheet.xlsx — renumber monster/character IDs to start at 0
# (instead of 1) on all three sheets, and refresh the selection/active-sheet
# UI state to match.

$wb = $excel.ActiveWorkbook

$sheetNames = @("OneStarDatas", "TwoStarDatas", "ThreeStarDatas")

# 1) Decrement the ID column (A2:A30) by 1 on every sheet so IDs run 0..28
#    instead of 1..29. Only column A changes; everything else is untouched.
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 30; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}

# 2) Update each sheet's selection to A2:A30 (anchored at A2), matching the
#    new ID column that was just edited.
$wsThree = $wb.Worksheets.Item("ThreeStarDatas")
$wsThree.Activate()
$wsThree.Range("A2:A30").Select()

$wsTwo = $wb.Worksheets.Item("TwoStarDatas")
$wsTwo.Activate()
$wsTwo.Range("A2:A30").Select()

$wsOne = $wb.Worksheets.Item("OneStarDatas")
$wsOne.Activate()
$wsOne.Range("A2:A30").Select()

# Activating OneStarDatas last makes it both the active sheet (tabSelected /
# activeTab) and leaves its selection as A2:A30 — matching the target state.
